# Applies the Thu Sep  7 23:52:38 UTC 2023 "Updated cryptos list" data
# refresh: updates Price/Volume(1h) figures for rows 2-45 and
# refreshes + reorders the last six coin rows (46-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (D) and Volume(1h) (E) columns for rows 2-45 ---
$ws.Range("D2").Value = "26.065.92"
$ws.Range("E2").Value = "  +1.02%  "

$ws.Range("D3").Value = "1.638.25"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("D4").Value = "'0.995"
$ws.Range("E4").Value = "  -0.59%  "

$ws.Range("D5").Value = "'215.64"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "'0.505"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "'0.0634"
$ws.Range("E9").Value = "  -0.64%  "

$ws.Range("D10").Value = "'19.83"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D11").Value = "'0.0789"
$ws.Range("E11").Value = "  -0.52%  "

$ws.Range("E12").Value = "  -0.35%  "

$ws.Range("D13").Value = "1.863.77"
$ws.Range("E13").Value = "  +0.16%  "

$ws.Range("D14").Value = "1.626.88"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").Value = "'0.553"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").Value = "0.0₃0764"
$ws.Range("E16").Value = "  -0.58%  "

$ws.Range("D17").Value = "'63.28"
$ws.Range("E17").Value = "  +0.56%  "

$ws.Range("D18").Value = "26.035.57"
$ws.Range("E18").Value = "  +0.82%  "

$ws.Range("E19").Value = "  -0.65%  "

$ws.Range("D20").Value = "'4.45"
$ws.Range("E20").Value = "  -0.43%  "

$ws.Range("D21").Value = "'193.24"
$ws.Range("E21").Value = "  -0.43%  "

$ws.Range("D22").Value = "'10.03"
$ws.Range("E22").Value = "  +0.66%  "

$ws.Range("E23").Value = "  +0.90%  "

$ws.Range("D24").Value = "'0.995"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("E25").Value = "  -2.07%  "

$ws.Range("D26").Value = "'142.61"
$ws.Range("E26").Value = "  +0.07%  "

$ws.Range("E27").Value = "  +0.48%  "

$ws.Range("D28").Value = "'6.91"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("D29").Value = "'15.57"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("E31").Value = "  +1.08%  "

$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "  -0.61%  "

$ws.Range("E33").Value = "  -0.28%  "

$ws.Range("E34").Value = "  +0.72%  "

$ws.Range("D35").Value = "'2.40"
$ws.Range("E35").Value = "  +0.42%  "

$ws.Range("D36").Value = "'0.907"
$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("D37").Value = "1.140.62"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").Value = "'0.550"
$ws.Range("E38").Value = "  +0.55%  "

$ws.Range("D39").Value = "'2.50"
$ws.Range("E39").Value = "  -1.42%  "

$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("D41").Value = "'0.995"
$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("E42").Value = "  +0.20%  "

$ws.Range("D43").Value = "'100.55"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").Value = "'0.793"
$ws.Range("E44").Value = "  -1.88%  "

$ws.Range("D45").Value = "1.773.67"
$ws.Range("E45").Value = "  +0.16%  "

# --- Rows 46-51: coins reshuffled with refreshed Price/Volume(1h) ---
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'55.89"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -7.83%  "

$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.47"
$ws.Range("E48").Value = "  +5.25%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0516"
$ws.Range("E49").Value = "  +2.13%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.417"
$ws.Range("E50").Value = "  +0.05%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.67"
$ws.Range("E51").Value = "  +1.78%  "
